$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; existing header (row1) becomes row2,
# existing data (row2) becomes row3.
$ws.Rows.Item(1).Insert()

# New label above the first (existing) 20-round ROC-AUC table.
$ws.Range("A1").Value = "H-2-Kd"

# Label for the new (second) table, directly below the first table's data row.
$ws.Range("A4").Value = "H-2-Kb"

# Re-create the header row (1-round .. 20-round) for the new table by
# copying the existing header row.
$ws.Range("A2:T2").Copy()
$ws.Range("A5").PasteSpecial()

# Write the new pearson-correlation data row for H-2-Kb.
$values = @(0.80439628137903896, 0.80698228033022501, 0.80998924469505496, 0.80771853238971703, 0.80949232742986799, 0.80749432802905796, 0.81200344057695595, 0.80146176962877402, 0.80798236409357804, 0.80706128250817399, 0.81322693708600402, 0.80572021388874904, 0.81479722236245999, 0.816846977161093, 0.81164607310565295, 0.80948464482669602, 0.80946574356911105, 0.81247017338392502, 0.809821259631885, 0.81240935974595296)

$col = 1
foreach ($v in $values) {
  $ws.Cells.Item(6, $col).Value = $v
  $col = $col + 1
}

# Match the highlighted-column formatting (E column) used in the first table.
$ws.Range("E3").Copy()
$ws.Range("E6").PasteSpecial(-4122)

# Restore selection similar to the authored workbook.
$ws.Range("I8").Select()
